# Update "想去人数" (want-to-go headcount) figures across sheets, matching
# the refreshed data output committed at 456a3b4.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1002
$ws.Range("F4").Value = 242
$ws.Range("F7").Value = 954
$ws.Range("F8").Value = 293
$ws.Range("F9").Value = 66
$ws.Range("F11").Value = 905
$ws.Range("F12").Value = 328
$ws.Range("F14").Value = 533
$ws.Range("F15").Value = 1381
$ws.Range("F17").Value = 1285
$ws.Range("F18").Value = 2945
$ws.Range("F19").Value = 322
$ws.Range("F22").Value = 762
$ws.Range("F24").Value = 1310
$ws.Range("F26").Value = 1081
$ws.Range("F28").Value = 3344

# ---- Sheet "演出" ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F14").Value = 3
$ws.Range("F15").Value = 10

# ---- Sheet "全部类型" ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 1002
$ws.Range("F7").Value = 242
$ws.Range("F11").Value = 954
$ws.Range("F12").Value = 293
$ws.Range("F14").Value = 66
$ws.Range("F23").Value = 905
$ws.Range("F24").Value = 328
$ws.Range("F26").Value = 533
$ws.Range("F27").Value = 1381
$ws.Range("F29").Value = 1285
$ws.Range("F30").Value = 2945
$ws.Range("F31").Value = 322
$ws.Range("F34").Value = 762
$ws.Range("F36").Value = 1310
$ws.Range("F40").Value = 1081
$ws.Range("F42").Value = 3344
$ws.Range("F46").Value = 3
$ws.Range("F47").Value = 10
